$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename option columns in row 2 (new shared strings get created in this order)
$ws.Range("E2").Value = "option_02"
$ws.Range("F2").Value = "option_03"

# New header cell for "wrongOption" info column
$ws.Range("G1").Value = "n"
$ws.Range("G2").Value = "wrongOption"

# Move the trailing comment string from G to H and populate G with
# the wrongOption index (0-based) for each quiz row
$ws.Range("H3").Value = "//sound도 index랑 동일하게 "
$ws.Range("G3").Value = 0
$ws.Range("G4").Value = 1
$ws.Range("G5").Value = 2

# Update the active selection to match the authored state
[void]$ws.Range("K8").Select()
